# 07/01/25 Commit - added tc 30 software restrictions
# TC_TM_030_apply_security_network_SoftwareRestriction
#
# Adds a new "SoftwareRestriction" worksheet (right before "Sheet1", i.e.
# immediately after "ProxySettings") populated with the Installation &
# Uninstallation / Application / Browser restriction test-data rows, and
# makes it the active sheet. Also tweaks the ProxySettings sheet selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new worksheet right before "Sheet1" (after ProxySettings)
# ---------------------------------------------------------------------
$proxySettings = $wb.Worksheets.Item("ProxySettings")
$ws = $wb.Worksheets.Add($null, $proxySettings)
$ws.Name = "SoftwareRestriction"

# ---------------------------------------------------------------------
# 2) Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Template Name"
$ws.Range("B1").Value = "Select Tab"
$ws.Range("C1").Value = "Installation & Uninstallation Restriction"
$ws.Range("D1").Value = "Software Restriction"
$ws.Range("E1").Value = "Software Application  Name"
$ws.Range("F1").Value = "Browser Name"
$ws.Range("G1").Value = "Restriction Type"

# ---------------------------------------------------------------------
# 3) Data rows
# ---------------------------------------------------------------------
$ws.Range("A2").Formula = "=MasterTemplate"
$ws.Range("B2").Value = "Installation & Uninstallation"
$ws.Range("C2").Value = "Y"
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "NA"

$ws.Range("A3").Formula = "=MasterTemplate"
$ws.Range("B3").Value = "Application Restriction"
$ws.Range("C3").Value = "NA"
$ws.Range("D3").Value = "Allow All"
$ws.Range("E3").Value = "chrome.exe"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "NA"

$ws.Range("A4").Formula = "=MasterTemplate"
$ws.Range("B4").Value = "Application Restriction"
$ws.Range("C4").Value = "NA"
$ws.Range("D4").Value = "Deny All"
$ws.Range("E4").Value = "chrome.exe"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "NA"

$ws.Range("A5").Formula = "=MasterTemplate"
$ws.Range("B5").Value = "Browser Restriction"
$ws.Range("C5").Value = "NA"
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "Google Chrome"
$ws.Range("G5").Value = "Block all downloads"

# ---------------------------------------------------------------------
# 4) Formatting - header row (yellow fill, bordered) and body (bordered)
# ---------------------------------------------------------------------
$headerLeft = $ws.Range("A1:C1")
$headerLeft.Interior.Color = 65535
$headerLeft.Borders.Item(7).LineStyle = 1
$headerLeft.Borders.Item(10).LineStyle = 1
$headerLeft.Borders.Item(8).LineStyle = 1

$headerRight = $ws.Range("D1:G1")
$headerRight.Interior.Color = 65535
$headerRight.Borders.Item(7).LineStyle = 1
$headerRight.Borders.Item(10).LineStyle = 1

$body = $ws.Range("A2:G5")
$body.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 5) Column widths (best-fit approximation of the authored widths)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.15
$ws.Columns.Item(2).ColumnWidth = 25.5
$ws.Columns.Item(3).ColumnWidth = 35.65
$ws.Columns.Item(4).ColumnWidth = 18.5
$ws.Columns.Item(5).ColumnWidth = 25.5
$ws.Columns.Item(6).ColumnWidth = 14.0
$ws.Columns.Item(7).ColumnWidth = 17.8

# ---------------------------------------------------------------------
# 6) Sheet view: make it the active/selected sheet, zoomed to 130%,
#    with D3 selected.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("D3").Select()

# ---------------------------------------------------------------------
# 7) ProxySettings selection changes to a block range (no single active
#    cell), matching the diff.
# ---------------------------------------------------------------------
$proxySettings.Activate()
$proxySettings.Range("A1:D2").Select()

# Re-activate the new sheet so it stays the active/selected tab.
$ws.Activate()
